# Generate Report for handoff
# Adds a new "handoff" row (757b0ebf-bf45-461b-b0bb-aadff4df32e6.md) above the
# pre-existing ".localization-config" row on every sheet, and fills in the
# handoff details (handoff file, handoff datetime, handoff reason) on the
# per-locale sheets for that new file.

$wb = $excel.ActiveWorkbook

$commit = "3ab82e35a75efda0fd9f9acfcb3c170377871c00"
$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/$commit"

$newFile = "757b0ebf-bf45-461b-b0bb-aadff4df32e6.md"
$oldFile = ".localization-config"

$zhHandoffFile = "757b0ebf-bf45-461b-b0bb-aadff4df32e6.85adad9b00b33310dbdaf6673550a79ffb99affb.zh-cn.xlf"
$deHandoffFile = "757b0ebf-bf45-461b-b0bb-aadff4df32e6.85adad9b00b33310dbdaf6673550a79ffb99affb.de-de.xlf"

$zhHandoffDate = "2016-01-07 08:01:06"
$deHandoffDate = "2016-01-07 08:01:18"

$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Drop the existing hyperlink on A2 so it can be rebuilt pointing at the
# new file, then move the old ".localization-config" entry down to row 3.
$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A3").Value = $oldFile
$ws.Range("B3").Value = "Not localized"
$ws.Range("C3").Value = "Not localized"

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("C2").Value = "Not yet handed off"
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$newFile", "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$oldFile", "", "", $oldFile)

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A3").Value = $oldFile
$ws.Range("B3").Value = "Not localized"
$ws.Range("D3").Value = $epoch
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Ignored"

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("C2").Value = $zhHandoffFile
$ws.Range("D2").Value = $zhHandoffDate
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$newFile", "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/$zhHandoffFile", "", "", $zhHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$oldFile", "", "", $oldFile)

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A3").Value = $oldFile
$ws.Range("B3").Value = "Not localized"
$ws.Range("D3").Value = $epoch
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Ignored"

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "Not yet handed off"
$ws.Range("C2").Value = $deHandoffFile
$ws.Range("D2").Value = $deHandoffDate
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$newFile", "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/$deHandoffFile", "", "", $deHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$oldFile", "", "", $oldFile)
